# gaza_nutrition_parameters.xlsx -- add two new parameter rows ("wt_loss_child_min",
# "wt_loss_child_max") plus a "pop" row to the "general" sheet, and flip which
# sheet/cell is the active selection (general tab becomes the active tab instead
# of scenarios).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")

# --- Row 6: wt_loss_child_min -------------------------------------------------
# Column A first (parameter name), then C (range), then B (shared description),
# then D (yes/no flag) -- this mirrors how the new rows were actually typed and
# keeps the shared-string table ordering lined up.
$ws.Range("A6").Value = "wt_loss_child_min"
$ws.Range("A7").Value = "wt_loss_child_max"

$ws.Range("C6").Value = "minimum of range"
$ws.Range("C7").Value = "maximum of range"

$ws.Range("B6").Value = "weight loss of children, relative to adults"
$ws.Range("B7").Value = "weight loss of children, relative to adults"

$ws.Range("D6").Value = "N"
$ws.Range("D7").Value = "N"

$ws.Range("E6").NumberFormat = "General"
$ws.Range("E6").Value = 0.7

$ws.Range("E7").NumberFormat = "General"
$ws.Range("E7").Value = 1.3

# --- Row 8: pop ----------------------------------------------------------------
$ws.Range("A8").Value = "pop"
$ws.Range("B8").Value = "total population of Gaza"
$ws.Range("C8").Value = "as of 7 Oct 2023"
$ws.Range("D8").Value = "N"

$ws.Range("E8").NumberFormat = "#,##0"
$ws.Range("E8").Value = 2226544

# --- Column widths (best-fit for the newly widened label/description columns) --
$ws.Columns.Item(1).ColumnWidth = 17
$ws.Columns.Item(2).ColumnWidth = 36.833333333333336
$ws.Columns.Item(3).ColumnWidth = 16.833333333333332
$ws.Columns.Item(5).ColumnWidth = 13.666666666666666

# --- Switch the active tab/selection from "scenarios" to "general" -------------
$wsScenarios = $wb.Worksheets.Item("scenarios")
$wsScenarios.Range("B14").Select()

$ws.Activate()
$ws.Range("A9").Select()
